# Updated symbol list on Mon Jan 23 21:52:52 UTC 2023 with GitHub Actions
#
# This script updates the "Price" (column D) and "Volume(1h)" (column E)
# figures for the crypto ranking sheet. The source data is stored as
# plain text (e.g. "305.75", "1.14%") rather than native numbers/percents
# -- each value keeps whatever fixed number of decimal places the scraper
# produced (e.g. "5.060", "0.0001300"), which would be lost if Excel
# re-parsed it as a real number or percentage. To avoid Excel's automatic
# "looks like a number" / "looks like a percentage" reinterpretation, each
# cell is briefly switched to a Text ("@") number format while the literal
# string is written, then ClearFormats() restores the cell's original
# (default/general) formatting -- so only the displayed text changes, not
# the cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "305.75" },
    @{ Cell = "E2"; Value = "1.14%" },
    @{ Cell = "D3"; Value = "36.22" },
    @{ Cell = "E3"; Value = "-2.06%" },
    @{ Cell = "D4"; Value = "5.060" },
    @{ Cell = "E4"; Value = "1.68%" },
    @{ Cell = "D5"; Value = "0.07942" },
    @{ Cell = "E5"; Value = "2.92%" },
    @{ Cell = "D6"; Value = "2.237" },
    @{ Cell = "E6"; Value = "7.26%" },
    @{ Cell = "D7"; Value = "8.003" },
    @{ Cell = "E7"; Value = "0.85%" },
    @{ Cell = "D8"; Value = "0.9288" },
    @{ Cell = "E8"; Value = "1.57%" },
    @{ Cell = "D9"; Value = "0.09824" },
    @{ Cell = "E9"; Value = "2.30%" },
    @{ Cell = "D10"; Value = "0.1882" },
    @{ Cell = "E10"; Value = "2.13%" },
    @{ Cell = "D11"; Value = "0.09081" },
    @{ Cell = "E11"; Value = "6.73%" },
    @{ Cell = "D12"; Value = "0.03710" },
    @{ Cell = "E12"; Value = "5.76%" },
    @{ Cell = "D13"; Value = "0.09922" },
    @{ Cell = "E13"; Value = "-0.40%" },
    @{ Cell = "D14"; Value = "0.001434" },
    @{ Cell = "E14"; Value = "-2.52%" },
    @{ Cell = "D15"; Value = "0.005601" },
    @{ Cell = "E15"; Value = "-1.74%" },
    @{ Cell = "D16"; Value = "3.449" },
    @{ Cell = "E16"; Value = "-0.53%" },
    @{ Cell = "D17"; Value = "4.144" },
    @{ Cell = "E17"; Value = "3.25%" },
    @{ Cell = "E18"; Value = "18.60%" },
    @{ Cell = "E19"; Value = "-0.35%" },
    @{ Cell = "D20"; Value = "0.1317" },
    @{ Cell = "E20"; Value = "-0.64%" },
    @{ Cell = "D21"; Value = "5.108" },
    @{ Cell = "E21"; Value = "7.34%" },
    @{ Cell = "E22"; Value = "2.19%" },
    @{ Cell = "D23"; Value = "0.04543" },
    @{ Cell = "E23"; Value = "-0.99%" },
    @{ Cell = "D24"; Value = "0.001237" },
    @{ Cell = "E24"; Value = "0.44%" },
    @{ Cell = "D25"; Value = "0.004784" },
    @{ Cell = "E25"; Value = "-6.28%" },
    @{ Cell = "D26"; Value = "0.0001301" },
    @{ Cell = "E26"; Value = "-7.29%" },
    @{ Cell = "D39"; Value = "0.01923" },
    @{ Cell = "E39"; Value = "9.16%" },
    @{ Cell = "D40"; Value = "0.04937" },
    @{ Cell = "E40"; Value = "7.28%" },
    @{ Cell = "D41"; Value = "0.007839" },
    @{ Cell = "E41"; Value = "5.09%" },
    @{ Cell = "D42"; Value = "0.1396" },
    @{ Cell = "E42"; Value = "0.32%" },
    @{ Cell = "D43"; Value = "0.007797" },
    @{ Cell = "E43"; Value = "0.80%" },
    @{ Cell = "D44"; Value = "0.002171" },
    @{ Cell = "E44"; Value = "0.28%" },
    @{ Cell = "D45"; Value = "0.01143" },
    @{ Cell = "E45"; Value = "10.56%" },
    @{ Cell = "D46"; Value = "0.00006255" },
    @{ Cell = "E46"; Value = "-1.07%" },
    @{ Cell = "D47"; Value = "0.00000000749" },
    @{ Cell = "D48"; Value = "51.85" },
    @{ Cell = "E48"; Value = "37.66%" },
    @{ Cell = "D49"; Value = "0.001798" },
    @{ Cell = "E49"; Value = "-10.08%" },
    @{ Cell = "D50"; Value = "0.00002098" },
    @{ Cell = "D51"; Value = "0.0001998" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}

